# Update the "Förändrad" (last changed) date column (C) for every data row
# (rows 2-135) on the "Avverkningsanmälningar" sheet.
# The date serial value 45172 (2023-09-03) becomes 45175 (2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C135").Value = 45175
